$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2125.2903
$ws.Range("I15").Value = 2125.2903
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 6375.8709
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -6206.8709

$ws.Range("H76").Value = 3060
$ws.Range("I76").Value = 3023.4482
$ws.Range("J76").Value = 3325
$ws.Range("K76").Value = 3023.4482
$ws.Range("L76").Value = 3325
$ws.Range("M76").Value = -2708.4482
$ws.Range("N76").Value = -3955

$ws.Range("H79").Value = 3060
$ws.Range("I79").Value = 3023.4482
$ws.Range("J79").Value = 3325
$ws.Range("K79").Value = 3023.4482
$ws.Range("L79").Value = 3325
$ws.Range("M79").Value = -1931.4482
$ws.Range("N79").Value = -5509

$ws.Range("H137").Value = 1432.8235
$ws.Range("I137").Value = 1151.6428
$ws.Range("J137").Value = 2745
$ws.Range("K137").Value = 3454.9284
$ws.Range("L137").Value = 8235
$ws.Range("M137").Value = -904.9284000000002
$ws.Range("N137").Value = -13335

$ws.Range("H138").Value = 1850.59
$ws.Range("I138").Value = 788.4146
$ws.Range("J138").Value = 2588.712
$ws.Range("K138").Value = 2365.2438
$ws.Range("L138").Value = 7766.136
$ws.Range("M138").Value = 2774.7562
$ws.Range("N138").Value = -18046.136

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3949.8115
$ws.Range("I32").Value = 2947.6897
$ws.Range("J32").Value = 9233.728
$ws.Range("K32").Value = 2947.6897
$ws.Range("L32").Value = 9233.728
$ws.Range("M32").Value = -2660.6897
$ws.Range("N32").Value = -9807.728

$ws.Range("H48").Value = 100000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 100000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 100000
$ws.Range("N48").Value = -100768

$ws.Range("H61").Value = 1539.4565
$ws.Range("I61").Value = 1278.0303
$ws.Range("J61").Value = 2203.077
$ws.Range("K61").Value = 1278.0303
$ws.Range("L61").Value = 2203.077
$ws.Range("M61").Value = -1066.0303
$ws.Range("N61").Value = -2627.077

$ws.Range("H74").Value = 252303
$ws.Range("I74").Value = 252303
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 252303
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -251429
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 252303
$ws.Range("I77").Value = 252303
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 1261515
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -1257147
$ws.Range("N77").ClearContents()

$ws.Range("H122").Value = 1314.2549
$ws.Range("I122").Value = 953.5294
$ws.Range("J122").Value = 2035.7059
$ws.Range("K122").Value = 2860.5882
$ws.Range("L122").Value = 6107.1177
$ws.Range("M122").Value = -410.5882000000001
$ws.Range("N122").Value = -11007.1177

$ws.Range("H128").Value = 26906.863
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 26906.863
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 26906.863
$ws.Range("N128").Value = -36866.863

$ws.Range("H136").Value = 1539.4565
$ws.Range("I136").Value = 1278.0303
$ws.Range("J136").Value = 2203.077
$ws.Range("K136").Value = 3834.0909
$ws.Range("L136").Value = 6609.231000000001
$ws.Range("M136").Value = -1284.0909
$ws.Range("N136").Value = -11709.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 1000
$ws.Range("I128").Value = 1000
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 3000
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -510

$ws.Range("H134").Value = 12999.167
$ws.Range("I134").Value = 15998.333
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 47994.999
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -45459.999
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22729012
$ws.Range("I31").Value = 40001184
$ws.Range("J31").Value = 2468.8948
$ws.Range("K31").Value = 40001184
$ws.Range("L31").Value = 2468.8948
$ws.Range("M31").Value = -40000889
$ws.Range("N31").Value = -3058.8948

$ws.Range("H34").Value = 22729012
$ws.Range("I34").Value = 40001184
$ws.Range("J34").Value = 2468.8948
$ws.Range("K34").Value = 40001184
$ws.Range("L34").Value = 2468.8948
$ws.Range("M34").Value = -40000982
$ws.Range("N34").Value = -2872.8948

$ws.Range("H132").Value = 3055.96
$ws.Range("I132").Value = 2444.5557
$ws.Range("J132").Value = 4628.143
$ws.Range("K132").Value = 7333.6671
$ws.Range("L132").Value = 13884.429
$ws.Range("M132").Value = -4803.6671
$ws.Range("N132").Value = -18944.429

$ws.Range("H134").Value = 1942.6666
$ws.Range("I134").Value = 1976.4667
$ws.Range("J134").Value = 1604.6666
$ws.Range("K134").Value = 5929.4001
$ws.Range("L134").Value = 4813.9998
$ws.Range("M134").Value = -3394.4001
$ws.Range("N134").Value = -9883.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 9000
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 9000
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 27000
$ws.Range("N95").Value = -31118

$ws.Range("H131").Value = 872.25
$ws.Range("J131").Value = 908.1573
$ws.Range("L131").Value = 2724.4719
$ws.Range("N131").Value = -12804.4719

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2255.95
$ws.Range("I102").Value = 2200.375
$ws.Range("J102").Value = 2478.25
$ws.Range("K102").Value = 2200.375
$ws.Range("L102").Value = 2478.25
$ws.Range("M102").Value = -578.375
$ws.Range("N102").Value = -5722.25

$ws.Range("H122").Value = 29622.805
$ws.Range("I122").Value = 38693.15
$ws.Range("J122").Value = 2411.7778
$ws.Range("K122").Value = 116079.45
$ws.Range("L122").Value = 7235.3334
$ws.Range("M122").Value = -113629.45
$ws.Range("N122").Value = -12135.3334

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H132").Value = 2237.7297
$ws.Range("I132").Value = 2199.4688
$ws.Range("J132").Value = 2482.6
$ws.Range("K132").Value = 6598.4064
$ws.Range("L132").Value = 7447.799999999999
$ws.Range("M132").Value = -4068.4064
$ws.Range("N132").Value = -12507.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4046.8333
$ws.Range("I7").Value = 3959.7273
$ws.Range("J7").Value = 5005
$ws.Range("K7").Value = 3959.7273
$ws.Range("L7").Value = 5005
$ws.Range("M7").Value = -3847.7273
$ws.Range("N7").Value = -5229

$ws.Range("H16").Value = 932.6667
$ws.Range("I16").Value = 932.6667
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 932.6667
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -762.6667
$ws.Range("N16").ClearContents()

$ws.Range("H40").Value = 2115.8
$ws.Range("I40").Value = 1919.75
$ws.Range("J40").Value = 2900
$ws.Range("K40").Value = 1919.75
$ws.Range("L40").Value = 2900
$ws.Range("M40").Value = -1783.75
$ws.Range("N40").Value = -3172

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H82").Value = 1129.4445
$ws.Range("I82").Value = 1338.8182
$ws.Range("J82").Value = 985.5
$ws.Range("K82").Value = 1338.8182
$ws.Range("L82").Value = 985.5
$ws.Range("M82").Value = -977.8181999999999
$ws.Range("N82").Value = -1707.5

$ws.Range("H85").Value = 1129.4445
$ws.Range("I85").Value = 1338.8182
$ws.Range("J85").Value = 985.5
$ws.Range("K85").Value = 1338.8182
$ws.Range("L85").Value = 985.5
$ws.Range("M85").Value = -90.81819999999993
$ws.Range("N85").Value = -3481.5

$ws.Range("H122").Value = 5454.5293
$ws.Range("I122").Value = 6185.1665
$ws.Range("J122").Value = 3701
$ws.Range("K122").Value = 18555.4995
$ws.Range("L122").Value = 11103
$ws.Range("M122").Value = -16105.4995
$ws.Range("N122").Value = -16003

$ws.Range("H123").Value = 20275
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 20275
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 20275
$ws.Range("N123").Value = -30075

$ws.Range("H126").Value = 4046.8333
$ws.Range("I126").Value = 3959.7273
$ws.Range("J126").Value = 5005
$ws.Range("K126").Value = 11879.1819
$ws.Range("L126").Value = 15015
$ws.Range("M126").Value = -9409.1819
$ws.Range("N126").Value = -19955

$ws.Range("H130").Value = 20750
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 20750
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 20750
$ws.Range("N130").Value = -30790

$ws.Range("H132").Value = 6629.727
$ws.Range("I132").Value = 10001
$ws.Range("J132").Value = 4703.2856
$ws.Range("K132").Value = 30003
$ws.Range("L132").Value = 14109.8568
$ws.Range("M132").Value = -27473
$ws.Range("N132").Value = -19169.8568

$ws.Range("H136").Value = 17545718
$ws.Range("I136").Value = 23811018
$ws.Range("J136").Value = 2881
$ws.Range("K136").Value = 71433054
$ws.Range("L136").Value = 8643
$ws.Range("M136").Value = -71430504
$ws.Range("N136").Value = -13743

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1052.2222
$ws.Range("I81").Value = 1072.5
$ws.Range("J81").Value = 890
$ws.Range("K81").Value = 2145
$ws.Range("L81").Value = 1780
$ws.Range("M81").Value = -1084
$ws.Range("N81").Value = -3902

$ws.Range("H84").Value = 1052.2222
$ws.Range("I84").Value = 1072.5
$ws.Range("J84").Value = 890
$ws.Range("K84").Value = 10725
$ws.Range("L84").Value = 8900
$ws.Range("M84").Value = -5421
$ws.Range("N84").Value = -19508

$ws.Range("H122").Value = 83733336
$ws.Range("I122").Value = 250000000
$ws.Range("J122").Value = 600002.5
$ws.Range("K122").Value = 750000000
$ws.Range("L122").Value = 1800007.5
$ws.Range("M122").Value = -749997550
$ws.Range("N122").Value = -1804907.5

$ws.Range("H128").Value = 30417
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 30417
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 30417
$ws.Range("N128").Value = -40377

$ws.Range("H132").Value = 4084466.5
$ws.Range("I132").Value = 5558460
$ws.Range("J132").Value = 2638.6924
$ws.Range("K132").Value = 16675380
$ws.Range("L132").Value = 7916.0772
$ws.Range("M132").Value = -16672850
$ws.Range("N132").Value = -12976.0772
